$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the instructional placeholder row (old row 2); this shifts the
# sample data row (old row 3) up to become row 2.
$ws.Rows("2:2").Delete()

# Update header row text to combine the previous instructional text into
# the header cells themselves.
$ws.Range("A1").Value = "상호"
$ws.Range("B1").Value = "대표이사"
$ws.Range("C1").Value = "법인등록번호`n000000-0000000"
$ws.Range("D1").Value = "설립년월일`nYYYY-MM-DD"
$ws.Range("E1").Value = "본점소재지`n서울 = 1`n경기 = 2`n충청북도 = 3`n충청남도 = 4`n강원도 = 5`n경상북도 = 6`n경상남도 = 7`n전라북도 = 8`n전라남도 = 9`n인천 = 10`n세종 = 11`n대전 = 12`n대구 = 13`n울산 = 14`n광주 = 15`n부산 = 16`n제주 = 17"
$ws.Range("F1").Value = "상세주소`n(법인등기부등본상)"
$ws.Range("G1").Value = "자본금(백만원)"
$ws.Range("H1").Value = "금융상품유형`n대출 = 1`n시설대여 및 연불판매 = 2`n할부 = 3`n어음할인 = 4`n매출채권 매입 = 5`n지급보증 = 6`n기타 대출성상품 = 7"
$ws.Range("I1").Value = "계약일자`nYYYY-MM-DD"
$ws.Range("J1").Value = "위탁예정기간`nYYYY-MM-DD"

# Header cells that now contain multi-line instructional text need
# word-wrap turned on to match their new taller row.
$ws.Range("C1:E1").WrapText = $true
$ws.Range("H1:J1").WrapText = $true

# Row 1 grows much taller to fit the long multi-line instructional text.
$ws.Rows("1:1").RowHeight = 313.2

# Update sample data row (now row 2) - company name, CEO name, address
$ws.Range("A2").Value = "베이직"
$ws.Range("B2").Value = "홍길동"
$ws.Range("F2").Value = "용산구 한남동 221-14"

# Move the active selection as the author left it after editing.
[void]$ws.Range("E4").Select()

Write-Host "done"
